$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Zebronics - Wireless Bluetooth Headset"
$ws.Range("B1").Value = "Rs. 847"

$ws.Range("A2").Value = "TUNE AUDIO U&I NICE Neckband Wireless With Mic Headphones/Earphones Beige and Black"
$ws.Range("B2").Value = "Rs. 749"

$ws.Range("A3").Value = "NBOX STAR WIRELESS NECKBAND WITH DOLBY EFFECT BASS SOUND IPX5 WITH MASSIVE MUSIC PLAYBACK WITH 1 YEAR WARRANTY BLUETOOTH HEADPHONE,BLUETOOTH EARPHONE,BLUETOOTH NECKBAND"
$ws.Range("B3").Value = "Rs. 749"

$ws.Range("A4").Value = "NBOX INVICTUS Neckband Wireless With Mic Headphones/Earphones Gold"
$ws.Range("B4").Value = "Rs. 999"

$ws.Range("A5").Value = "boAt Airdopes 131/138 On Ear True Wireless (TWS) 15 Hours Playback IPX7(Water Resistant) Active Noise cancellation -Bluetooth V 5.0 Black"
$ws.Range("B5").Value = "Rs. 1,299"
